# "added entropy grid visualization" -- extend the 2-column entropy grid
# on the sheet with two more rows of data, and scroll/select to show the
# newly added area, matching how Excel records the view after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended right after the existing grid (row 33 was the
# previous last row), growing the used range from A1:B33 to A1:B35.
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = 5
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = 4

# Reflect the post-edit viewport/selection: scrolled down so row 7 is the
# first visible row, with the active cell/selection on C36 (just past the
# new data), matching the sheetView recorded for the edited workbook.
$ws.Range("C36").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
